$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Расчетно-пояснительная записка на ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

$d.Content.Find.Execute(
    "91",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "96", 2
)
